# Auto-generated COM-interop edit script for LOB1224.docx
$d = $word.ActiveDocument
$vtab = [char]11

# 1) Paragraph 6 (Normal): PT objectives text -> PT 'Programa resumido' text
$d.Paragraphs.Item(6).Range.Text = 'Elementos de teoria e história do planejamento urbano. Teoria e prática do planejamento ambiental; Planejamento ambiental como indutor de desenvolvimento sustentável; Aplicações da teoria do planejamento a problemas ambientais e urbanos; Legislação e política ambiental urbana.'

# 2) Paragraph 7 (Normal, italic): EN objectives text -> EN 'Programa resumido' text
$d.Paragraphs.Item(7).Range.Text = 'Theory elements and history of urban planning. Theory and practice of environmental planning; environmental planning as an inducer of sustainable development; environmental theory applied to urban problems; legislation and urban environmental policy.'

# 3) Paragraph 9 (ListBullet, under 'Docente(s) Responsavel(eis)'): teacher name -> PT objectives text
$d.Paragraphs.Item(9).Range.Text = 'Propiciar ao discente uma visão integrada do processo de planejamento com um enfoque ambiental aplicado no urbanismo; apresentar e discutir conceitos, projetos práticos e metodologias relacionadas às etapas e fases do planejamento ambiental e à gestão ambiental urbana; apresentar os instrumentos do planejamento, gestão e política ambiental urbana'

# 4) Paragraph 11 (Normal, under 'Programa resumido'): PT resumido text -> PT 'Programa' text
$d.Paragraphs.Item(11).Range.Text = 'Introdução ao planejamento e gestão ambiental. Origens da teoria e prática do planejamento. Natureza do planejamento e suas relações com a geografia, política, economia, sociedade, cultura e meio ambiente. Análises, estudos e proposições relativas às diversas formas de crescimento e expansão urbanas; Elementos para estruturação ambiental da cidade; Etapas, estruturas e instrumentos do planejamento ambiental; Indicadores ambientais e planejamento; Participação pública no planejamento ambiental; Política Nacional do Meio Ambiente (Lei n°6938/1981); Sistema Nacional de Unidades de Conservação (Lei n°9985/2000); Estatuto da Cidade (Lei n°10.257/2001); Zoneamento Ambiental; EIA e EIV como instrumentos inovadores; Novos conceitos e princípios de planos diretores urbano-ambientais;'

# 5) Paragraph 12 (Normal, italic, under 'Programa resumido'): EN resumido text -> EN objectives text
$d.Paragraphs.Item(12).Range.Text = 'Provide an integrated comprehension about planning process from an environmental approach applied for urban planning; introduce and discuss concepts, practical projects and methodologies related to stages and phases of environmental planning and urban environmental management; introducing instruments of planning, management and urban environmental policy.'

# 6) Paragraph 14 (Normal, under 'Programa'): PT programa text -> Metodo text (2 text runs + line break)
$d.Paragraphs.Item(14).Range.Text = 'Aulas teóricas e práticas, visitas técnicas e exercícios dirigidos. ' + $vtab + 'Avaliação baseada em provas, exercícios e trabalhos práticos e relatórios.'

# 7) Paragraph 17 (ListBullet, 'Avaliacao'): restructure the 3 plain-text runs that
#    follow the bold 'Metodo:'/'Criterio:'/'Norma de recuperacao:' labels. The bold
#    label runs themselves are unchanged and used as stable anchors.
$p17 = $d.Paragraphs.Item(17)

# 7a) Span after 'Norma de recuperacao: ' (was 'Provas e/ou...', no trailing break)
#     -> becomes the full bibliography text (no trailing break)
$rng = $p17.Range
[void]$rng.Find.Execute('Norma de recuperação: ')
$rng.Collapse(0)
$spanStart = $rng.Start
$spanEnd = $p17.Range.End
$target = $d.Range($spanStart, $spanEnd)
$target.Text = 'Bibliografia básica:' + $vtab + 'AGRA FILHO, S,S. Planejamento e Gestão Ambiental no Brasil. Os Instrumentos da Política Nacional do Meio Ambiente, Rio de Janeiro, Elsevier, 2014' + $vtab + 'FRANCO, M.A.R., Planejamento ambiental para a cidade sustentável, Ed. Annablume, 2000' + $vtab + 'DEAK, C., SHIFFER, S.T.R., O processo de urbanização no Brasil, EDUSP, 1999' + $vtab + 'IBGE, Instituto Brasileiro de Geografia e Estatística. Indicadores de Desenvolvimento Sustentável. Rio de Janeiro, IBGE, 2012.' + $vtab + 'MOTA, S., Urbanização e meio ambiente, ABES Associação Brasileira de Engenharia Sanitária, 1999' + $vtab + 'MENEZES, C.L., Desenvolvimento urbano e meio ambiente, Papirus, 1996' + $vtab + 'PHILLIPI, Jr.A; MALHEIROS, T.F. Indicadores de Sustentabilidade e Gestão Ambiental. Editora Manole, 2012.' + $vtab + 'SANTOS, M. A Urbanização Brasileira. 3 ed. São Paulo: HUCITEC, 1993. 155p' + $vtab + 'SANTOS, R.F., Planejamento ambiental: teoria e prática, Editora Oficina de textos, 2004' + $vtab + 'SECCHI, L. Análise de Políticas Públicas. Diagnóstico de Problemas, Recomendações de Soluções., São Paulo, Cengage Learning, 2016' + $vtab + 'SOUZA, M.L. Mudar a Cidade: Uma introdução crítica ao planejamento e à gestão urbanos. Rio de Janeiro, Bertrand Brasil, 2003.' + $vtab + 'VILLAÇA, F. Uma contribuição para a história do planejamento urbano no Brasil. In: DEAK, C; SCHIFFER, S.R (org) O processo de urbanização no Brasil. São Paulo, EDUSP, 1999.' + $vtab + '' + $vtab + 'Bibliografia complementar:' + $vtab + 'ALLEN, A., YOU, N., Sustainable urbanization – bridging the green and brown agendas, DPU, University College London, 2002' + $vtab + 'ACSELRAD, H., Conflitos ambientais no Brasil, Fundação Henrich Boll, 2004' + $vtab + 'BARDET, G., O urbanismo, Papirus, 1990' + $vtab + 'BUARQUE, S.C., LIMA, R.R.A.; Manual de estratégia de desenvolvimento para aglomerações urbanas, Brasília, IPEA, 2005' + $vtab + 'MENEGAT, R; ALMEIDA, G. Desenvolvimento Sustentável e Gestão Ambiental nas Cidades. Porto Alegre, Editora UFRGS, 2004.'

# 7b) Span between 'Criterio: ' and 'Norma de recuperacao: ' (was 'Media ponderada...<br>')
#     -> becomes 'Provas e/ou exercicios dirigidos.' + break
$rng = $p17.Range
[void]$rng.Find.Execute('Critério: ')
$rng.Collapse(0)
$spanStart = $rng.Start
$rng2 = $d.Range($spanStart, $p17.Range.End)
[void]$rng2.Find.Execute('Norma de recuperação: ')
$spanEnd = $rng2.Start
$target = $d.Range($spanStart, $spanEnd)
$target.Text = 'Provas e/ou exercícios dirigidos.' + $vtab

# 7c) Span between 'Metodo: ' and 'Criterio: ' (was 'Aulas...<br>Avaliacao...<br>')
#     -> becomes 'Media ponderada...' + break
$rng = $p17.Range
[void]$rng.Find.Execute('Método: ')
$rng.Collapse(0)
$spanStart = $rng.Start
$rng2 = $d.Range($spanStart, $p17.Range.End)
[void]$rng2.Find.Execute('Critério: ')
$spanEnd = $rng2.Start
$target = $d.Range($spanStart, $spanEnd)
$target.Text = 'Média ponderada das notas atribuídas às provas, exercícios e trabalhos práticos e relatórios.' + $vtab

# 8) Paragraph 19 (Normal, under 'Bibliografia'): bibliography text -> teacher name
$d.Paragraphs.Item(19).Range.Text = '9146830 - Danúbia Caporusso Bargos'

Write-Output 'done'
